$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.521.24"
$ws.Range("E2").Value = "  -4.38%  "

# Row 3
$ws.Range("D3").Value = "3.258.19"
$ws.Range("E3").Value = "  -5.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.69"
$ws.Range("E5").Value = "  -2.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.23"
$ws.Range("E6").Value = "  -4.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  -2.20%  "

# Row 9
$ws.Range("D9").Value = "3.251.70"
$ws.Range("E9").Value = "  -5.22%  "

# Row 10
$ws.Range("E10").Value = "  -8.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.585"
$ws.Range("E11").Value = "  -4.41%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.21"
$ws.Range("E12").Value = "  -6.94%  "

# Row 13
$ws.Range("E13").Value = "  -6.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "634.56"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.55"
$ws.Range("E15").Value = "  -5.20%  "

# Row 16
$ws.Range("D16").Value = "3.779.92"
$ws.Range("E16").Value = "  -5.02%  "

# Row 17
$ws.Range("D17").Value = "65.407.47"

# Row 18
$ws.Range("E18").Value = "  -3.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.66"
$ws.Range("E19").Value = "  -1.86%  "

# Row 20
$ws.Range("D20").Value = "3.255.26"
$ws.Range("E20").Value = "  -5.51%  "

# Row 21
$ws.Range("E21").Value = "  -6.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.900"
$ws.Range("E22").Value = "  -3.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.78"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.66"
$ws.Range("E24").Value = "  +7.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.93"
$ws.Range("E25").Value = "  -7.38%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").Value = "  -6.36%  "

# Row 27
$ws.Range("E27").Value = "  -5.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  -2.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.66"
$ws.Range("E29").Value = "  -4.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.29"
$ws.Range("E30").Value = "  -5.67%  "

# Row 31
$ws.Range("E31").Value = "  -2.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  -5.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.99"
$ws.Range("E33").Value = "  -4.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "543.53"
$ws.Range("E34").Value = "  +8.84%  "

# Row 35
$ws.Range("E35").Value = "  -3.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.89"
$ws.Range("E37").Value = "  -5.86%  "

# Row 38
$ws.Range("D38").Value = "3.591.66"
$ws.Range("E38").Value = "  -1.04%  "

# Row 39
$ws.Range("E39").Value = "  -1.58%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  -4.48%  "

# Row 41
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0715"
$ws.Range("E41").Value = "  -8.04%  "

# Row 42
$ws.Range("E42").Value = "  -1.83%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.26"
$ws.Range("E43").Value = "  -5.91%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "31.88"
$ws.Range("E44").Value = "  -6.11%  "

# Row 45
$ws.Range("E45").Value = "  -0.36%  "

# Row 46
$ws.Range("E46").Value = "  -8.47%  "

# Row 47
$ws.Range("E47").Value = "  -4.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  -3.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  -6.38%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("E51").Value = "  +1.52%  "
